$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 8 (SVM / VGG): hyperparameters, Train Accuracy, CV Accuracy change; Test Accuracy stays the same ---
$ws.Range("C8").Value = "C: 10.0, class_weight: balanced, gamma: scale, kernel: rbf"
$ws.Range("D8").Value = 1
$ws.Range("E8").Value = 0.8697293921731891

# --- Insert a new row above row 12 (copy row 13's formatting) so the XGBoost/VGG row
#     (and everything below it) shifts down by one ---
$ws.Rows.Item(13).Copy()
$ws.Rows.Item(12).Insert()

# Row 12 becomes "XGBoost / Top3 Features" reusing the same tuned hyperparameters text as the VGG row
$ws.Range("A12").Value = "XGBoost"
$ws.Range("B12").Value = "Top3 Features"
$ws.Range("C12").Value = "learning_rate: 0.5, max_depth: 3, min_child_weight: 1, n_estimators: 300"
$ws.Range("D12").Value = 1
$ws.Range("E12").Value = 0.8597356369691923
$ws.Range("F12").Value = 0.8488918737407656

# Row 13 keeps being "XGBoost / VGG" (unchanged values, just shifted down from old row 12)
$ws.Range("A13").Value = "XGBoost"
$ws.Range("B13").Value = "VGG"
$ws.Range("C13").Value = "learning_rate: 0.5, max_depth: 3, min_child_weight: 1, n_estimators: 300"
$ws.Range("D13").Value = 1
$ws.Range("E13").Value = 0.8637334165972801
$ws.Range("F13").Value = 0.8542646071188718

# --- Append a new row 17: Random Forest / Top3 Features (copy row 16's formatting) ---
$ws.Rows.Item(16).Copy()
$ws.Rows.Item(17).Insert()
$ws.Range("A17").Value = "Random Forest"
$ws.Range("B17").Value = "Top3 Features"
$ws.Range("C17").Value = "bootstrap: False, max_depth: 80, max_features: sqrt, min_samples_leaf: 2, n_estimators: 500"
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 0.7932697751873439
$ws.Range("F17").Value = 0.7575554063129617
